$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "word 예시" sheet had its English-example column (C) and Korean-
# translation column (D) authored in the wrong order. Swap the two whole
# columns (values + column widths) by cutting column C and re-inserting
# it after column D - a real physical column swap, not a per-cell value
# copy, so the column widths move with their data instead of being
# recomputed/rounded.
$ws.Columns.Item(3).Cut()
$ws.Columns.Item(5).Insert()

# Restore the active selection Excel leaves behind after this kind of
# edit.
$ws.Range("D8").Select()
